$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Chamber diameter" label to "Chamber radius" for both tables.
$ws.Range("C13").Value = "Chamber radius"
$ws.Range("F13").Value = "Chamber radius"

# Move the active selection to match the recorded cursor position after edit.
$ws.Range("J10").Select()
